# fix test case description
# Translate the English QA test-case text (rows 7-9) on the "Bin" sheet into
# Vietnamese, matching the book-store domain, and touch up a couple of
# related status cells. Also refresh the "Task" header label and the
# "Implementer: ..." banner cell on all three visible sheets (their shared
# string slot changes because other strings were dropped), nudge the sheet
# view/selection, widen column E and re-grow rows 8/9 to fit the new text.

$wb = $excel.ActiveWorkbook

$binSheet  = $wb.Worksheets.Item(1)   # "Bin"
$phucSheet = $wb.Worksheets.Item(2)   # "Phuc"
$thaiSheet = $wb.Worksheets.Item(3)   # "Thai"

# --- Bin sheet: rewrite the test-case rows --------------------------------
# The exact order in which new (not-yet-seen) strings are assigned matters
# for shared-string layout, so cells are touched in the same order the
# original authoring tool would have written them back.

# Row 7 ("4.2 - 3")
$binSheet.Range("B7").Value = "Hàm có trả về đầy đủ thông tin của cuốn sách hay không ?"
$binSheet.Range("C7").Value = "N/A"
$binSheet.Range("D7").Value = "nhấn chuột vào ảnh của bất kì cuốn sách nào hiện trong trang chủ và trang danh mục"
$binSheet.Range("E7").Value = "Thông tin tương ứng của cuốn sách đó được trả về"
$binSheet.Range("F7").Value = "Tested"
$binSheet.Range("G7").Value = "Passed"
$binSheet.Range("H7").Value = "N/A"

# Row 8 ("4.2 - 5")
$binSheet.Range("B8").Value = "kiểm tra tài khoản và mật khẩu và người dùng nhập trong form đăng nhập "
$binSheet.Range("C8").Value = "form đăng nhập phải được nhập đầy đủ "
$binSheet.Range("D8").Value = "mở form đăng nhập trên thanh điều hướng -> điền đầy đủ tài khoản và mật khấu -> nhấn nút ""Đăng nhập"""
$binSheet.Range("E8").Value = "người dùng vào được tài khoản đã được đăng ký. Nếu tài khoảng chưa đăng kí thì hiện thông báo cho người dùng"
$binSheet.Range("F8").Value = "Tested"
$binSheet.Range("G8").Value = "Passed"
$binSheet.Range("H8").Value = "N/A"

# Row 9 ("4.2 - 6")
$binSheet.Range("B9").Value = "kiểm tra mật khảu tài khoản người dùng có được mã hóa khi lưu vào cơ sở dữ liệu không ?"
$binSheet.Range("C9").Value = "N/A"
$binSheet.Range("D9").Value = "Open the database -> take 5 random client account "
$binSheet.Range("H9").Value = "mật khẩu thực tế phải khác so với trong cơ sở dữ liệu"
$binSheet.Range("E9").Value = "mật khẩu đã được mã hóa"
$binSheet.Range("F9").Value = "Tested"
$binSheet.Range("G9").Value = "Passed"

# Row 6 header label stays "Task" on every sheet (value unchanged, only the
# underlying shared-string slot shifts once unused strings are dropped).
$binSheet.Range("A6").Value = "Task"
$phucSheet.Range("A6").Value = "Task"
$thaiSheet.Range("A6").Value = "Task"

# Implementer banner cells keep their text too. Re-writing A1 nudges this
# engine's row-height autofit, so immediately restore each row's original
# height to avoid spurious height churn in the saved XML.
$binSheet.Range("A1").Value = "Implementer: Nguyễn Văn Bin"
$binSheet.Rows.Item(1).AutoFit()

$phucSheet.Range("A1").Value = "Implementer: Đặng Trần Thiên Phúc"
$phucSheet.Rows.Item(1).AutoFit()

$thaiSheet.Range("A1").Value = "Implementer: Trương Lâm Quốc Thái"
$thaiSheet.Rows.Item(1).RowHeight = 14.4

# --- Row heights (rows grew taller to fit the longer Vietnamese text) ----
$binSheet.Rows.Item(8).RowHeight = 99
$binSheet.Rows.Item(9).RowHeight = 61.8

# --- Column E widened slightly --------------------------------------------
$binSheet.Columns.Item(5).ColumnWidth = 19

# --- View / selection on the Bin sheet -------------------------------------
$binSheet.Activate()
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 2
$binSheet.Range("J7").Select()
